$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

$ws.Range("B54").Value = 713
$ws.Range("C54").Value = 275
$ws.Range("D54").Value = 836
$ws.Range("E54").Value = 319
